# Updates cryptos list: Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.708.36"
$ws.Range("E2").Value = "  +3.15%  "
$ws.Range("D3").Value = "1.788.87"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'222.74"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'32.49"
$ws.Range("E8").Value = "  +8.60%  "
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").Value = "'0.0685"
$ws.Range("E10").Value = "  +3.65%  "
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("D12").Value = "2.045.00"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").Value = "'11.00"
$ws.Range("E13").Value = "  +10.58%  "
$ws.Range("D14").Value = "1.767.26"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "34.706.78"
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").Value = "'4.29"
$ws.Range("E17").Value = "  +3.19%  "
$ws.Range("D18").Value = "'68.48"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "'252.80"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").Value = "0.0₃0783"
$ws.Range("E20").Value = "  +6.69%  "
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "'10.47"
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").Value = "'158.88"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'7.04"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "'0.0515"
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").Value = "'3.74"
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("D35").Value = "1.431.20"
$ws.Range("E35").Value = "  -2.93%  "
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").Value = "'82.79"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").Value = "'2.81"
$ws.Range("E40").Value = "  +4.20%  "
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").Value = "'0.901"
$ws.Range("E42").Value = "  +2.40%  "
$ws.Range("D43").Value = "'2.05"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("E46").Value = "  +4.44%  "
$ws.Range("D47").Value = "1.941.68"
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("D48").Value = "'104.23"
$ws.Range("E48").Value = "  +7.57%  "
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "'11.96"
$ws.Range("E50").Value = "  +2.34%  "
$ws.Range("D51").Value = "'49.78"
$ws.Range("E51").Value = "  -1.66%  "
